$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.387.97"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "3.590.24"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.86"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("E6").Value = "  +18.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "653.17"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.431"
$ws.Range("E8").Value = "  +7.62%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.06"
$ws.Range("E10").Value = "  +4.33%  "
$ws.Range("D11").Value = "3.587.28"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.37"
$ws.Range("E12").Value = "  +4.72%  "
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.50"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "4.256.72"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "97.120.48"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").Value = "3.593.26"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.17"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.536"
$ws.Range("E22").Value = "  +10.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "519.72"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  +5.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "104.05"
$ws.Range("E27").Value = "  +9.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.18"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.177"
$ws.Range("E29").Value = "  +22.18%  "
$ws.Range("D30").Value = "3.782.20"
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.03"
$ws.Range("E32").Value = "  +4.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.188"
$ws.Range("E34").Value = "  +6.17%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.02"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.581"
$ws.Range("E37").Value = "  +3.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "616.96"
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("E39").Value = "  +1.25%  "
$ws.Range("E40").Value = "  -4.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.154"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.933"
$ws.Range("E43").Value = "  +2.66%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.12"
$ws.Range("E45").Value = "  +6.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.440"
$ws.Range("E46").Value = "  +41.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0447"
$ws.Range("E47").Value = "  +7.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  +1.30%  "
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.68"
$ws.Range("E50").Value = "  +5.76%  "
$ws.Range("E51").Value = "  +7.76%  "
